# "Most grades are complete" - fill in newly-graded assignment scores.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - max-points header row: assignments 13-16 (O:R) now have point values, S stays 0
$ws.Range("O2").Value = 25
$ws.Range("P2").Value = 25
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0

# Row 3 - Arevalo, Andres
$ws.Range("J3").Value = 20
$ws.Range("N3").Value = 22
$ws.Range("O3").Value = 25
$ws.Range("P3").Value = 25

# Row 5
$ws.Range("L5").Value = 14
$ws.Range("N5").Value = 22
$ws.Range("P5").Value = 13

# Row 9
$ws.Range("O9").Value = 25
$ws.Range("P9").Value = 25

# Row 10
$ws.Range("L10").Value = 20
$ws.Range("N10").Value = 22
$ws.Range("P10").Value = 25

# Row 12
$ws.Range("N12").Value = 22
$ws.Range("O12").Value = 25
$ws.Range("P12").Value = 25

# Row 13
$ws.Range("G13").Value = 13
$ws.Range("H13").Value = 10
$ws.Range("I13").Value = 20
$ws.Range("J13").Value = 20
$ws.Range("N13").Value = 22

# Row 14
$ws.Range("D14").Value = 41
$ws.Range("E14").Value = 19
$ws.Range("F14").Value = 24
$ws.Range("G14").Value = 13
$ws.Range("H14").Value = 10
$ws.Range("I14").Value = 20
$ws.Range("K14").Value = 4

# Row 18
$ws.Range("P18").Value = 25

# Row 19
$ws.Range("D19").Value = 41
$ws.Range("I19").Value = 20
$ws.Range("K19").Value = 4
$ws.Range("L19").Value = 15

# Row 25
$ws.Range("O25").Value = 25

# Row 28
$ws.Range("D28").Value = 41
$ws.Range("E28").Value = 19
$ws.Range("F28").Value = 24
$ws.Range("G28").Value = 13
$ws.Range("H28").Value = 10
$ws.Range("I28").Value = 20
$ws.Range("J28").Value = 20
$ws.Range("K28").Value = 4
$ws.Range("L28").Value = 22
$ws.Range("N28").Value = 22
$ws.Range("O28").Value = 25
$ws.Range("P28").Value = 25

# Row 29
$ws.Range("I29").Value = 20
$ws.Range("N29").Value = 22
$ws.Range("O29").Value = 25
$ws.Range("P29").Value = 25

# Row 30
$ws.Range("H30").Value = 10
$ws.Range("I30").Value = 20
$ws.Range("L30").Value = 22
$ws.Range("N30").Value = 22
$ws.Range("O30").Value = 25
$ws.Range("P30").Value = 25

# Row 38
$ws.Range("I38").Value = 20

# Row 43
$ws.Range("G43").Value = 13
$ws.Range("H43").Value = 10
$ws.Range("I43").Value = 20
$ws.Range("J43").Value = 20
$ws.Range("L43").Value = 22
$ws.Range("N43").Value = 22
$ws.Range("O43").Value = 25
$ws.Range("P43").Value = 25

# View state: zoomed to 85% with Q30 as the active selection
$excel.ActiveWindow.Zoom = 85
$excel.Goto($ws.Range("Q30"))
